$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (e.g. "565.50", "0.591") are written
# as literal text rather than being auto-coerced to numbers by Excel, so
# formatting like trailing zeros / fixed decimal places is preserved exactly
# as it appears in the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.429.34"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.333.71"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "190.21"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").Value = "565.50"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "3.326.43"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.591"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "48.04"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "3.870.62"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "607.28"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "66.491.47"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "18.13"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "3.335.06"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "11.17"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "18.81"
$ws.Range("E23").Value = "  +11.56%  "
$ws.Range("D24").Value = "5.20"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "9.77"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +8.70%  "
$ws.Range("D32").Value = "4.04"
$ws.Range("E32").Value = "  +6.23%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.18"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "562.60"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "57.47"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.720.72"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "0.0₃0734"
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "34.16"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.131"
$ws.Range("E41").Value = "  +4.97%  "
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D43").Value = "3.43"
$ws.Range("E43").Value = "  +7.74%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.71"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  +4.76%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "2.61"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  +3.90%  "

# Restore the default (unstyled) look for column D now that the text values
# are committed, so no lingering text-number-format style remains on the cells.
$ws.Range("D2:D51").Style = "Normal"
